$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark from the "Conceptual database
#    design:" paragraph (it will be re-added at the very end of the
#    document, after the new paragraph we add below).
# -------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# -------------------------------------------------------------------
# 2) Merge the "<w:tab/>" run and the "Doctor:" text run that live in
#    the same paragraph into a single run, keeping the tab character
#    as a real <w:tab/> element.
# -------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$found = $false
$n = $d.Paragraphs.Count
$targetText = [char]9 + "Doctor:" + [char]13
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ((-not $found) -and $t -eq $targetText) {
        $found = $true
        $pPrXml = '<w:pPr ' + $wNs + '><w:spacing w:line="240" w:lineRule="auto"/><w:contextualSpacing/></w:pPr>'
        $runXml = '<w:r ' + $wNs + '><w:tab/><w:t>Doctor:</w:t></w:r>'
        $paraXml = '<w:p ' + $wNs + '>' + $pPrXml + $runXml + '</w:p>'
        $p.Range.InsertXML($paraXml)
    }
}

# -------------------------------------------------------------------
# 3) Append a new paragraph describing the 4th query after the existing
#    empty "ListParagraph" paragraph that follows "Query Description",
#    then move the "_GoBack" bookmark to the end of that new paragraph.
# -------------------------------------------------------------------
$newPara = $d.Paragraphs.Add()

$run1 = "Our query chose to figure out which receptionists scheduled appointments for patients receiving Viagra prescribed by a given pharmacist. This allows a pharmacist to discover who is aware of the medication a patient is taking. Without this feature, it may be difficult to discover "
$run2 = "who is involved in the stages of a patient" + [char]8217 + "s hospital lifetime in the case of an emergency. Starting with a Pharmacist ID, it checks the "
$run3 = "PrescriptionRecord"
$run4 = " table for a matching ID and Viagra as the prescription. Then it checks the Patient ID associated with that prescription, and finds appointments for that patient, checking the receptionist that made the appointment and printing them out."

$newParaXml = '<w:p ' + $wNs + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">' + $run1 + '</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">' + $run2 + '</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>' + $run3 + '</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">' + $run4 + '</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$newPara.Range.InsertXML($newParaXml)
